$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 34, shifting existing rows 34-66 down to 35-67
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record
$ws.Cells.Item(34,1).Value = 9
$ws.Cells.Item(34,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(34,3).Value = "Metropolitana"
$ws.Cells.Item(34,4).Value = 45072
$ws.Cells.Item(34,5).Value = 13
$ws.Cells.Item(34,6).Value = "Fruta"
$ws.Cells.Item(34,7).Value = 100107
$ws.Cells.Item(34,8).Value = "Otros"
$ws.Cells.Item(34,9).Value = 100107001
$ws.Cells.Item(34,10).Value = "Caqui"
$ws.Cells.Item(34,11).Value = "Fuyu"
$ws.Cells.Item(34,12).Value = "Primera"
$ws.Cells.Item(34,13).Value = 470
$ws.Cells.Item(34,14).Value = 13000
$ws.Cells.Item(34,15).Value = 14000
$ws.Cells.Item(34,16).Value = 13468
$ws.Cells.Item(34,17).Value = "`$/caja 16 kilos granel"
$ws.Cells.Item(34,18).Value = "Región de O'Higgins"
$ws.Cells.Item(34,19).Value = 842
$ws.Cells.Item(34,20).Value = 16
